$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix mislabeled model name in first table (row 6 was a duplicate of row 5's
# "fixed + spatial/time/ar1" label; rename it to the correct model name).
$ws.Range("A6").Value = "fixed + spatial/time"

# Update second table title (D1) from "pefa" to "rlha"
$ws.Range("D1").Value = "rlha"

# Row 5 of second table: change label and value
#   D5: "fixed + spatial/time/ar1" -> "fixed + spatial/time/iid"
#   E5: 188.01 -> 187.74 (value that used to live in the now-removed row 6)
$ws.Range("D5").Value = "fixed + spatial/time/iid"
$ws.Range("E5").Value = 187.74

# Remove the old row 6 entries for the second table (D6/E6), since that
# content got folded into row 5 above.
$ws.Range("D6:E6").Clear()

# New formula cells
$ws.Range("F4").Formula = "=E3-E5"
$ws.Range("B8").Formula = "=B6-B3"

# Update the selection shown in the saved view
$ws.Range("B9").Select()
